$d = $word.ActiveDocument

# --- 1. Split "Geometry - Normal/bump map information" into two runs:
#        "Geometry - Normal/bump map" + "/Opacity"
$found = $d.Content.Find.Execute("Geometry – Normal/bump map information", $true, $false, $false, $false, $false, $true, 1, $false, "Geometry – Normal/bump map", 2)

$geoPara = $d.Paragraphs(10).Range
$insertPos = $geoPara.End - 1
$insertPoint = $d.Range($insertPos, $insertPos)
$insertPoint.InsertAfter("/Opacity")

# Force the newly-typed text into its own run (distinct from the preceding
# run) by toggling a character property that round-trips back to the
# paragraph's default formatting.
$newRunRange = $d.Range($insertPos, $insertPos + 8)
$newRunRange.Font.Bold = $true
$newRunRange.Font.Bold = $false

# --- 2. Insert three new list paragraphs after the "Bump mapping tab..." paragraph.
$bumpPara = $d.Paragraphs(11).Range
$bumpPara.InsertParagraphAfter()
$opacityPara = $d.Paragraphs(12).Range
$opacityPara.Text = "Opacity map is used for cutting out transparent parts of an object. (NOTE: Make sure to disable “opaque” under the object’s shape node!)"
$opacityPara.ListFormat.ListLevelNumber = 3

$opacityPara2 = $d.Paragraphs(12).Range
$opacityPara2.InsertParagraphAfter()
$transPara = $d.Paragraphs(13).Range
$transPara.Text = "Transmission – Refractive information"
$transPara.ListFormat.ListLevelNumber = 2

$transPara2 = $d.Paragraphs(13).Range
$transPara2.InsertParagraphAfter()
$refractPara = $d.Paragraphs(14).Range
$refractPara.Text = "Similar to opacity, but used for refractive objects like glass and water."
$refractPara.ListFormat.ListLevelNumber = 3

Write-Output "edit complete"
